$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: "time_taken", styled like the other header cells (copy format from E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:39:19.750318",
    "2021-10-05 13:39:19.750330",
    "2021-10-05 13:39:19.750333",
    "2021-10-05 13:39:19.750336",
    "2021-10-05 13:39:19.750339",
    "2021-10-05 13:39:19.750342",
    "2021-10-05 13:39:19.750345",
    "2021-10-05 13:39:19.750347",
    "2021-10-05 13:39:19.750350",
    "2021-10-05 13:39:19.750353",
    "2021-10-05 13:39:19.750355",
    "2021-10-05 13:39:19.750358",
    "2021-10-05 13:39:19.750361",
    "2021-10-05 13:39:19.750363",
    "2021-10-05 13:39:19.750366",
    "2021-10-05 13:39:19.750368",
    "2021-10-05 13:39:19.750371",
    "2021-10-05 13:39:19.750374",
    "2021-10-05 13:39:19.750376",
    "2021-10-05 13:39:19.750379",
    "2021-10-05 13:39:19.750382",
    "2021-10-05 13:39:19.750384",
    "2021-10-05 13:39:19.750387",
    "2021-10-05 13:39:19.750390",
    "2021-10-05 13:39:19.750392",
    "2021-10-05 13:39:19.750395",
    "2021-10-05 13:39:19.750397",
    "2021-10-05 13:39:19.750400",
    "2021-10-05 13:39:19.750402",
    "2021-10-05 13:39:19.750405",
    "2021-10-05 13:39:19.750408",
    "2021-10-05 13:39:19.750410",
    "2021-10-05 13:39:19.750414",
    "2021-10-05 13:39:19.750416",
    "2021-10-05 13:39:19.750419",
    "2021-10-05 13:39:19.750421",
    "2021-10-05 13:39:19.750424",
    "2021-10-05 13:39:19.750426",
    "2021-10-05 13:39:19.750429",
    "2021-10-05 13:39:19.750431",
    "2021-10-05 13:39:19.750434",
    "2021-10-05 13:39:19.750437",
    "2021-10-05 13:39:19.750440",
    "2021-10-05 13:39:19.750442",
    "2021-10-05 13:39:19.750445",
    "2021-10-05 13:39:19.750447",
    "2021-10-05 13:39:19.750450",
    "2021-10-05 13:39:19.750452",
    "2021-10-05 13:39:19.750455",
    "2021-10-05 13:39:19.750457",
    "2021-10-05 13:39:19.750460",
    "2021-10-05 13:39:19.750463",
    "2021-10-05 13:39:19.750466",
    "2021-10-05 13:39:19.750469",
    "2021-10-05 13:39:19.750472",
    "2021-10-05 13:39:19.750474",
    "2021-10-05 13:39:19.750477",
    "2021-10-05 13:39:19.750480",
    "2021-10-05 13:39:19.750482",
    "2021-10-05 13:39:19.750485",
    "2021-10-05 13:39:19.750487",
    "2021-10-05 13:39:19.750490",
    "2021-10-05 13:39:19.750493",
    "2021-10-05 13:39:19.750495",
    "2021-10-05 13:39:19.750499",
    "2021-10-05 13:39:19.750502",
    "2021-10-05 13:39:19.750505",
    "2021-10-05 13:39:19.750508"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
